$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = 'questions = [
    {
        "title": "The following elements are added to a queue (in order): 30, 70, 30, 80, 50. What would the queue contain after three delete operations?",
        "ques_type": 2,
        "options": [
            "80, 50",
            "30, 70",
            "NULL",
            "30, 70, 30"
        ],
        "score": "80, 50"
    },
    {
        "title": "For which of the following elements can a binary search tree not be constructed?",
        "ques_type": 2,
        "options": [
            "10, 20, 30, 12, 34, 23",
            "1.55, 2.33, 4.55, 9.12, 3.11, 2.15",
            "\u201cKanesha\u201d, \u201cDinesh\u201d, \u201cRaul\u201d, \u201cMark\u201d, \u201cMimi\u201d",
            "32, +, 35, /, 4, - 44"
        ],
        "score": "32, +, 35, /, 4, - 44"
    },
    {
        "title": "An algorithm needs n input parameters, 1 output parameter, 22 local variables, and an auxiliary array of n elements. Which of the following should be considered for evaluating its space complexity?",
        "ques_type": 2,
        "options": [
            "Input parameters",
            "Output parameters",
            "Local variables",
            "Array of n elements"
        ],
        "score": "Array of n elements"
    },
    {
        "title": "Which of the statements is correct about the following code snippet? Employee e = new Employee ( \"Samuel\", 24, \"Oakland Street, New Jersey\" )",
        "ques_type": 2,
        "options": [
            "An object of Employee class is created on heap.",
            "A copy constructor function of Employee class is called.",
            "An object of Employee class is created in stack.",
            "Variable e would contain the name, age, and address of the employee."
        ],
        "score": "An object of Employee class is created on heap."
    }
]'

# Row 2 (old shared-string cell) goes away; row 1's old bold/bordered/
# centered numeric cell is reset to the default style before it takes on
# the reformatted questions text as a (new) shared string.
$ws.Range("A2").ClearContents()
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $text
$ws.Rows.Item(1).AutoFit()
